$wb = $excel.ActiveWorkbook

# --- Sheet "Cardiac" ---
$cardiac = $wb.Worksheets.Item("Cardiac")

# Row 2 and Row 3: swap the B text (and A2 becomes TRUE)
$cardiac.Range("A2").Value = $true
$cardiac.Range("B2").Value = "Do you have any PMHx? (counts as 2 independent minor features)"
$cardiac.Range("B3").Value = "Pain not worse with exertion (requires they clarify exercise 1hr after meal)"

# Row 8: clear A8 (was TRUE, now blank)
$cardiac.Range("A8").ClearContents()

# Row 9: remove stray space before "(food gets stuck..."
$cardiac.Range("B9").Value = "Alternative cause of esoph dysphagia becomes obvious(food gets stuck or relieved by regurgitation of food)"

# Row 18: fix typo "soley" -> "solely"
$cardiac.Range("B18").Value = "Pain worse with exertion (without clarifying that it only occurs solely within an hour of eating)"

# --- Sheet "CREST" ---
$crest = $wb.Worksheets.Item("CREST")

# Row 5: clear A5 (was TRUE, now blank)
$crest.Range("A5").ClearContents()

# Row 6: set A6 to TRUE (was blank)
$crest.Range("A6").Value = $true

# Row 11: set A11 to TRUE (was blank)
$crest.Range("A11").Value = $true
